$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking price strings
# (e.g. '213.18') are not auto-converted into floating point numbers.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = '27.714.69'
$ws.Range("E2").Value = '  +1.06%  '

$ws.Range("D3").Value = '1.645.91'
$ws.Range("E3").Value = '  +0.19%  '

$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("D5").Value = '213.18'
$ws.Range("E5").Value = '  +0.56%  '

$ws.Range("E6").Value = '  -0.95%  '

$ws.Range("E7").Value = '  +0.11%  '

$ws.Range("D8").Value = '23.24'
$ws.Range("E8").Value = '  -0.17%  '

$ws.Range("E9").Value = '  +0.65%  '

$ws.Range("E10").Value = '  +0.61%  '

$ws.Range("D11").Value = '0.0893'
$ws.Range("E11").Value = '  +0.29%  '

$ws.Range("D12").Value = '1.878.71'
$ws.Range("E12").Value = '  +0.23%  '

$ws.Range("D13").Value = '1.641.44'
$ws.Range("E13").Value = '  +0.18%  '

$ws.Range("E14").Value = '  +0.25%  '

$ws.Range("D15").Value = '0.562'
$ws.Range("E15").Value = '  +0.69%  '

$ws.Range("D16").Value = '64.78'
$ws.Range("E16").Value = '  +0.80%  '

$ws.Range("D17").Value = '27.692.49'
$ws.Range("E17").Value = '  +1.07%  '

$ws.Range("D18").Value = '231.90'
$ws.Range("E18").Value = '  +1.44%  '

$ws.Range("D19").Value = '0.0₃0725'
$ws.Range("E19").Value = '  +1.01%  '

$ws.Range("D20").Value = '7.62'
$ws.Range("E20").Value = '  +1.46%  '

$ws.Range("E21").Value = '  +0.15%  '

$ws.Range("E22").Value = '  -0.28%  '

$ws.Range("D23").Value = '10.17'
$ws.Range("E23").Value = '  +9.03%  '

$ws.Range("E24").Value = '  -3.67%  '

$ws.Range("D25").Value = '150.17'
$ws.Range("E25").Value = '  +1.49%  '

$ws.Range("D26").Value = '6.93'
$ws.Range("E26").Value = '  -0.03%  '

$ws.Range("E27").Value = '  -2.55%  '

$ws.Range("D28").Value = '15.66'
$ws.Range("E28").Value = '  +0.78%  '

$ws.Range("E29").Value = '  +0.16%  '

$ws.Range("E30").Value = '  +0.85%  '

$ws.Range("E31").Value = '  -0.22%  '

$ws.Range("E32").Value = '  +1.14%  '

$ws.Range("D33").Value = '1.440.49'
$ws.Range("E33").Value = '  +2.30%  '

$ws.Range("E34").Value = '  +1.52%  '

$ws.Range("E35").Value = '  +1.81%  '

$ws.Range("E36").Value = '  -1.12%  '

$ws.Range("D37").Value = '0.571'
$ws.Range("E37").Value = '  +1.71%  '

$ws.Range("E38").Value = '  -0.01%  '

$ws.Range("E39").Value = '  +0.13%  '

$ws.Range("D40").Value = '0.884'
$ws.Range("E40").Value = '  +11.99%  '

$ws.Range("E41").Value = '  +0.81%  '

$ws.Range("D43").Value = '67.48'
$ws.Range("E43").Value = '  +4.41%  '

$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").Value = '5.59'
$ws.Range("E44").Value = '  +2.06%  '

$ws.Range("B45").Value = 'MXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D45").Value = '2.25'
$ws.Range("E45").Value = '  +1.58%  '

$ws.Range("B46").Value = 'RocketPoolETH'
$ws.Range("C46").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D46").Value = '1.787.83'
$ws.Range("E46").Value = '  +0.09%  '

$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").Value = '1.74'
$ws.Range("E47").Value = '  +5.97%  '

$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").Value = '0.0₆0108'
$ws.Range("E48").Value = '  +2.90%  '

$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").Value = '85.52'
$ws.Range("E49").Value = '  -2.18%  '

$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '7.82'
$ws.Range("E50").Value = '  +1.73%  '

$ws.Range("E51").Value = '  +0.35%  '

# Restore column D formatting to the default (General/no explicit style)
# now that the text values have been written.
$dRange.ClearFormats()

